# Functional Requirements.docx edit script
# Applies the changes described by the commit "Updated US and FR."

$d = $word.ActiveDocument

# wdLineSpaceMultiple constants used by Word for Paragraph.Format.LineSpacingRule
# 1 = wdLineSpace1pt5  -> serialises as <w:spacing w:line="360" w:lineRule="auto"/>
$wdLineSpace1pt5 = 1

# --- 1. Title paragraph: merge "Functional" + " " + "Requirements" into one run ---
$d.Content.Find.Execute("Functional Requirements", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Functional Requirements", 2) | Out-Null

# --- 2. Text edits (bullet list items) ---

# Bullet 1 (shopping cart item): add ", while shopping"
$d.Content.Find.Execute("which is consultable at any time.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "which is consultable at any time, while shopping.", 2) | Out-Null

# Bullet 3 (supermarkets list): "supermarkets which" -> "supermarkets, which"
$d.Content.Find.Execute("supermarkets which", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "supermarkets, which", 2) | Out-Null
# ... and "is sorted in order of shortest distance" -> "is sorted by distance"
$d.Content.Find.Execute("is sorted in order of shortest distance from a given location.", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "is sorted by distance from a given location.", 2) | Out-Null

# Bullet 4 (score): rewrite sentence
$d.Content.Find.Execute("The system shall provide a score based on product information to evaluate a shopping list’s health quality.", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The system shall provide a score for each cart, based on its products, to evaluate health quality.", 2) | Out-Null

# Bullet 5 (favourite products -> per-product score): rewrite sentence
$d.Content.Find.Execute("The system shall provide a list of favourite products, containing all products liked from the users.", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The system shall provide a score for each product, based on its label, to determine health quality. ", 2) | Out-Null

# Bullet 6 (calendar -> history of previous carts): rewrite sentence
$d.Content.Find.Execute("The system shall provide a calendar, in which users can insert all groceries they have planned to do in the future.", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The system shall provide a user’s history of previous carts, which is consultable for a new shop.", 2) | Out-Null

# --- 3. Line spacing: 1.5 lines on the blank paragraph after the title and on each bullet paragraph ---
$count = $d.Paragraphs.Count
for ($i = 2; $i -le $count; $i++) {
    $d.Paragraphs($i).Format.LineSpacingRule = $wdLineSpace1pt5
}

# --- 4. Remove the trailing empty paragraph (the one that had ind left=360) ---
# It is the final (now-empty) paragraph in the document; deleting the range
# that spans its preceding paragraph mark through the end of the story
# merges it away and leaves the previous paragraph's mark as the new end.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$prevEnd = $lastPara.Range.Start - 1
$docEnd = $d.Content.End
$d.Range($prevEnd, $docEnd).Delete() | Out-Null

Write-Output "edits applied"
